# Adds a new "2022" column (P) to the indicator table, mirroring the
# formatting of the existing "2021" column (O) but with a thousands-
# separator number format (#,##0.0 instead of 0.0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlRight = -4152, xlBottom = -4107
$xlRight = -4152
$xlBottom = -4107
$newFormat = "#,##0.0"

function Set-YearCell($sheet, $src, $tgt, $val) {
    $sheet.Range($src).Copy($sheet.Range($tgt))
    $sheet.Range($tgt).NumberFormat = $newFormat
    $sheet.Range($tgt).HorizontalAlignment = $xlRight
    $sheet.Range($tgt).VerticalAlignment = $xlBottom
    $sheet.Range($tgt).Value = $val
}

# Header year
$ws.Range("P4").Value = 2022

# Data rows (2022 figures for each indicator)
Set-YearCell $ws "O5"  "P5"  1188.7
Set-YearCell $ws "O6"  "P6"  263.9
Set-YearCell $ws "O7"  "P7"  263.2
Set-YearCell $ws "O8"  "P8"  12.4
Set-YearCell $ws "O9"  "P9"  "-"
Set-YearCell $ws "O10" "P10" 93
Set-YearCell $ws "O11" "P11" 171.5
Set-YearCell $ws "O12" "P12" 220.6
Set-YearCell $ws "O13" "P13" 159.3
Set-YearCell $ws "O14" "P14" 1.7
Set-YearCell $ws "O15" "P15" "-"
Set-YearCell $ws "O16" "P16" 3.1

# Match the author's final selection shown in the diff
[void]$ws.Range("Q7").Select()
